$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Grab references to the two existing sheets.
# ------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")

# ------------------------------------------------------------------
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE header on both sheets and
#    replace the full howstat.com scorecard URL with just the numeric
#    match code that used to be the MatchCode= query parameter.
# ------------------------------------------------------------------
$batting.Cells.Item(1, 4).Value = "MATCH_CODE"
for ($r = 2; $r -le 118; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $t = $cell.Text
    if ($t -match "MatchCode=(\d+)") {
        $cell.NumberFormat = "@"
        $cell.Value = $matches[1]
    }
}

$bowling.Cells.Item(1, 2).Value = "MATCH_CODE"
$bowlingLastRow = $bowling.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $t = $cell.Text
    if ($t -match "MatchCode=(\d+)") {
        $cell.NumberFormat = "@"
        $cell.Value = $matches[1]
    }
}

# ------------------------------------------------------------------
# 3. Insert the new "Player Info" sheet right before "ODI Batting"
#    so the tab order becomes:
#      Player Info, ODI Batting, ODI Bowling, ODI Batting Extra
# ------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

$playerInfo.Cells.Item(1, 1).Value = "ID"
$playerInfo.Cells.Item(1, 2).Value = "NAME"
$playerInfo.Cells.Item(1, 3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1, 4).Value = "BOWL_STYLE"

$hdr1 = $playerInfo.Range("A1:D1")
$hdr1.Font.Bold = $true
$hdr1.HorizontalAlignment = -4108
$hdr1.VerticalAlignment = -4160
$hdr1.Borders.LineStyle = 1

$playerInfo.Cells.Item(2, 1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "3513"
$playerInfo.Cells.Item(2, 2).Value = "Sarfaraz Ahmed"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Does Not Bowl | Unknown"

# ------------------------------------------------------------------
# 4. Append the new "ODI Batting Extra" sheet right after
#    "ODI Bowling" (fetch it fresh - its index shifted once the
#    Player Info sheet was inserted above).
# ------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowling)
$extra.Name = "ODI Batting Extra"

$extra.Cells.Item(1, 1).Value = "MATCH_CODE"
$extra.Cells.Item(1, 2).Value = "BATTING_POSITION"
$extra.Cells.Item(1, 3).Value = "NUM_4"
$extra.Cells.Item(1, 4).Value = "NUM_6"
$extra.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

$hdr2 = $extra.Range("A1:F1")
$hdr2.Font.Bold = $true
$hdr2.HorizontalAlignment = -4108
$hdr2.VerticalAlignment = -4160
$hdr2.Borders.LineStyle = 1

$extraRows = @(
    @("4227", "", "", "", "", "NO"),
    @("4237", 6, "0", "0", "0.37%", "NO"),
    @("4238", 8, "2", "0", "20.20%", "NO"),
    @("4241", 8, "1", "0", "1.89%", "NO"),
    @("4287", "", "", "", "", "NO"),
    @("4292", 6, "2", "0", "11.36%", "NO"),
    @("4294", "", "", "", "", "NO"),
    @("4297", 8, "3", "0", "6.18%", "NO"),
    @("4300", 5, "7", "2", "32.66%", "NO"),
    @("4304", 5, "1", "0", "7.62%", "NO"),
    @("4308", "", "", "", "", "NO"),
    @("4319", 5, "1", "0", "15.04%", "NO"),
    @("4324", 5, "0", "0", "5.66%", "NO"),
    @("4334", 8, "0", "0", "0.65%", "NO"),
    @("4337", "", "", "", "", "NO"),
    @("4340", 6, "1", "0", "7.83%", "NO"),
    @("4349", 7, "0", "0", "0.95%", "NO"),
    @("4375", 5, "0", "0", "2.62%", "NO"),
    @("4376", 4, "1", "0", "7.69%", "NO"),
    @("4460", 5, "2", "0", "4.06%", "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $cA = $extra.Cells.Item($r, 1)
    $cA.NumberFormat = "@"
    $cA.Value = $row[0]

    # BATTING_POSITION is a genuine number when present, blank otherwise.
    $cB = $extra.Cells.Item($r, 2)
    if ($row[1] -eq "") {
        $cB.Value = ""
    } else {
        $cB.Value = $row[1]
    }

    # NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL are stored as text.
    $cC = $extra.Cells.Item($r, 3)
    $cC.NumberFormat = "@"
    $cC.Value = $row[2]

    $cD = $extra.Cells.Item($r, 4)
    $cD.NumberFormat = "@"
    $cD.Value = $row[3]

    $cE = $extra.Cells.Item($r, 5)
    $cE.NumberFormat = "@"
    $cE.Value = $row[4]

    $cF = $extra.Cells.Item($r, 6)
    $cF.Value = $row[5]

    $r++
}

# ------------------------------------------------------------------
# 5. Restore the original active tab (first sheet) so the workbook
#    still opens on the same tab it did before the edit.
# ------------------------------------------------------------------
$playerInfo.Activate()
$playerInfo.Range("A1").Select() | Out-Null

